# "Data Semester Tahun 2025" sheet: split the old "Id Semester" column into
# two year columns ("Tahun 1" / "Tahun 2") and move the old "No" values into
# a new trailing "Semester" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 3) ---
$ws.Range("B3").Value = "Tahun 1"
$ws.Range("C3").Value = "Tahun 2"
$ws.Range("D3").Value = "Semester"

# Give the new header cell (D3) the same look as the other header cells.
$ws.Range("C3").Copy() | Out-Null
$ws.Range("D3").PasteSpecial(-4122) | Out-Null

# --- Data rows ---
$ws.Range("B4").Value = 2024
$ws.Range("C4").Value = 2025
$ws.Range("D4").Value = 1

$ws.Range("B5").Value = 2024
$ws.Range("C5").Value = 2025
$ws.Range("D5").Value = 2

# Give the new data cells (D4:D5) the same look as the other data cells.
$ws.Range("C4:C5").Copy() | Out-Null
$ws.Range("D4:D5").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# --- Column widths for the (now 3) data columns ---
$ws.Columns.Item(2).ColumnWidth = 9.283447
$ws.Columns.Item(3).ColumnWidth = 9.283447
$ws.Columns.Item(4).ColumnWidth = 10.568848

# --- Extend the title merge from A1:C1 to A1:D1 ---
$ws.Range("A1:C1").UnMerge()
$ws.Range("A1:D1").Merge()

# --- Selection matches the source workbook (last edited cell) ---
$ws.Range("D5").Select()
